$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused "need admin" comment from G35 and set the new actual hours.
$ws.Range("F20").Value = 1.5

$ws.Range("F35").Value = 2
$ws.Range("G35").Value = ""

$ws.Range("F38").Value = 1
$ws.Range("F39").Value = 0

$ws.Range("F42").Value = 2
$ws.Range("F43").Value = 3

$ws.Range("F46").Value = 10

# Restore the active selection as recorded after the edits.
$ws.Range("F21").Select()
